$d = $word.ActiveDocument

$lq = [char]0x201C
$rq = [char]0x201D

# 1. "Before a member can borrow a book, they must initiate a borrow book. "
#    -> "Before a member can borrow a book, they must initiate “Borrow Book”. "
# NOTE: the find span purposely starts at "they" so it swallows the
# _GoBack bookmark sitting right after "they" (matches Word's behaviour of
# dropping a bookmark that a replace operation's range passes through).
$replacement = "they must initiate " + $lq + "Borrow Book" + $rq + ". "
$found = $d.Content.Find.Execute("they must initiate a borrow book. ", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 1)

Write-Host "done find1: $found"

# 2. memberDAO -> memberMapDAO, loanDAO -> loanMapDAO, bookDAO -> bookMapDAO
$found = $d.Content.Find.Execute("memberDAO", $true, $false, $false, $false, $false, $true, 1, $false, "memberMapDAO", 1)
Write-Host "done find2 (memberDAO): $found"

$found = $d.Content.Find.Execute("loanDAO", $true, $false, $false, $false, $false, $true, 1, $false, "loanMapDAO", 1)
Write-Host "done find3 (loanDAO): $found"

$found = $d.Content.Find.Execute("bookDAO", $true, $false, $false, $false, $false, $true, 1, $false, "bookMapDAO", 1)
Write-Host "done find4 (bookDAO): $found"

# 3. Fix the typo "BorrowBokUI" -> "BorrowBookUI" for the occurrence followed
#    by " displayed" (no trailing period) only -- the other occurrence
#    (followed by " displayed.") stays misspelled per the source diff.
#    We insert the missing "o" in place (rather than a blanket Find/Replace)
#    so that the new _GoBack bookmark can be planted exactly where Word
#    left it (right after the newly-typed letter).
$r = $d.Content
$r.Find.Execute("BorrowBokUI displayed")
if (-not $r.Find.Found) {
    Write-Host "ERROR: BorrowBokUI displayed (no period) not found"
}
$wordStart = $r.Start
$insertPos = $wordStart + 8   # length of "BorrowBo"
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("o")
Write-Host "done find5 (BorrowBokUI -> BorrowBookUI)"

# Re-home the _GoBack bookmark at the point right after the inserted "o"
# (the old one was already destroyed above when the "they ... borrow book."
# span was replaced in step 1).
$bmPos = $insertPos + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Host "done bookmark relocate"
